$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 205, shifting existing rows 205:294 down to 206:295
$ws.Rows.Item(205).Insert()

# Populate the new row 205 with the new data point
$ws.Range("A205").Value = 10
$ws.Range("B205").Value = "Vega Modelo de Temuco"
$ws.Range("C205").Value = "La Araucanía"
$ws.Range("D205").Value = 44726
$ws.Range("E205").Value = 9
$ws.Range("F205").Value = 100112001
$ws.Range("G205").Value = "Berenjena"
$ws.Range("H205").Value = "Sin especificar"
$ws.Range("I205").Value = "Primera"
$ws.Range("J205").Value = 55
$ws.Range("K205").Value = 10000
$ws.Range("L205").Value = 10000
$ws.Range("M205").Value = 10000
$ws.Range("N205").Value = "$/caja 60 unidades"
$ws.Range("O205").Value = "Región de Arica y Parinacota"
$ws.Range("P205").Value = 167
$ws.Range("Q205").Value = 60
$ws.Range("R205").Value = "Hortaliza"
